$d = $word.ActiveDocument

# --- Part 1: remove the "HELLO WORLD" title paragraph entirely -------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.Delete()

# --- Part 2: merge the three runs "<link(", the URL, and ")Das ist" into a
# single run, while leaving the following " jetzt ein Link>" run untouched. --
$linkPara = $d.Paragraphs.Item(3)
$pStart = $linkPara.Range.Start

$mergeEnd = $pStart + 52

# Temporarily toggle Bold on the trailing run (" jetzt ein Link>") so that its
# formatting no longer matches the preceding runs; this stops the engine from
# coalescing it together with them when we edit the text just before it.
$tailRange = $d.Range($mergeEnd, $mergeEnd + 16)
$tailRange.Font.Bold = 1

# Insert a throwaway character right at the "<link(...)Das ist" / " jetzt..."
# boundary and remove it again. The insert+delete forces the engine to
# re-normalize runs around that boundary, merging the three identically
# formatted runs in front of it into one run (the differently formatted tail
# run is skipped because of the Bold toggle above).
$boundary = $d.Range($mergeEnd, $mergeEnd)
$boundary.InsertBefore("X")
$marker = $d.Range($mergeEnd, $mergeEnd + 1)
$marker.Delete()

# Restore the tail run's formatting back to normal (clears the Bold toggle
# without leaving any residue in the XML).
$tailRange2 = $d.Range($mergeEnd, $mergeEnd + 16)
$tailRange2.Font.Bold = 0

Write-Output $linkPara.Range.Text
